# Remove the 5 rows that correspond to stale/duplicate "Rope" (and one
# "Pruning saw") tool-count records, per the commit "added file exports
# for raw datasets". Rows are deleted in descending order so the row
# numbers of the as-yet-undeleted rows remain stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowsToDelete = @(66, 50, 31, 22, 20)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
